$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SBS Scores")

# Set column C to 0 for the rows that got a new bolus-dosage column value
$rows = @(4,6,7,9,10,11,12,14,15,16,17,18,19,20,21,22,23,25,26,27,28,31,32,34,35,38,40)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = 0
}

# Update the view: scroll back to top-left (A1) and move selection to G4
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G4").Select()
